$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.475.37'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '2.604.39'
$ws.Range("E3").Value = '  +6.08%  '

$ws.Range("E4").Value = '  +0.09%  '

$styleD5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.70'
$ws.Range("D5").Style = $styleD5
$ws.Range("E5").Value = '  +3.46%  '

$styleD6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.33'
$ws.Range("D6").Style = $styleD6
$ws.Range("E6").Value = '  +3.02%  '

$styleD7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.604'
$ws.Range("D7").Style = $styleD7
$ws.Range("E7").Value = '  +4.57%  '

$ws.Range("E8").Value = '  +0.12%  '

$styleD9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("D9").Style = $styleD9
$ws.Range("E9").Value = '  +12.15%  '

$styleD10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.42'
$ws.Range("D10").Style = $styleD10
$ws.Range("E10").Value = '  +10.91%  '

$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$styleD11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.68'
$ws.Range("D11").Style = $styleD11
$ws.Range("E11").Value = '  +1.87%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$styleD12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0847'
$ws.Range("D12").Style = $styleD12
$ws.Range("E12").Value = '  +7.30%  '

$styleD13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.17'
$ws.Range("D13").Style = $styleD13
$ws.Range("E13").Value = '  +13.52%  '

$ws.Range("D14").Value = '3.009.69'
$ws.Range("E14").Value = '  +7.11%  '

$ws.Range("E15").Value = '  +1.53%  '

$ws.Range("D16").Value = '2.614.64'
$ws.Range("E16").Value = '  +7.16%  '

$styleD17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.923'
$ws.Range("D17").Style = $styleD17
$ws.Range("E17").Value = '  +7.45%  '

$styleD18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '15.01'
$ws.Range("D18").Style = $styleD18
$ws.Range("E18").Value = '  +5.98%  '

$ws.Range("D19").Value = '46.626.11'
$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("E20").Value = '  +7.12%  '

$styleD21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.00'
$ws.Range("D21").Style = $styleD21
$ws.Range("E21").Value = '  +0.61%  '

$styleD22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.80'
$ws.Range("D22").Style = $styleD22
$ws.Range("E22").Value = '  +8.21%  '

$styleD23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '280.40'
$ws.Range("D23").Style = $styleD23
$ws.Range("E23").Value = '  +13.45%  '

$styleD24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.81'
$ws.Range("D24").Style = $styleD24
$ws.Range("E24").Value = '  +5.96%  '

$ws.Range("E25").Value = '  +8.15%  '

$styleD26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.17'
$ws.Range("D26").Style = $styleD26
$ws.Range("E26").Value = '  +10.08%  '

$styleD27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.89'
$ws.Range("D27").Style = $styleD27
$ws.Range("E27").Value = '  +33.91%  '

$styleD28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = $styleD28
$ws.Range("E28").Value = '  -0.21%  '

$styleD29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.04'
$ws.Range("D29").Style = $styleD29
$ws.Range("E29").Value = '  +0.64%  '

$styleD30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.66'
$ws.Range("D30").Style = $styleD30
$ws.Range("E30").Value = '  +8.19%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$styleD31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '39.14'
$ws.Range("D31").Style = $styleD31
$ws.Range("E31").Value = '  -1.84%  '

$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$styleD32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.25'
$ws.Range("D32").Style = $styleD32
$ws.Range("E32").Value = '  +1.18%  '

$styleD33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.35'
$ws.Range("D33").Style = $styleD33
$ws.Range("E33").Value = '  +12.67%  '

$styleD34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.68'
$ws.Range("D34").Style = $styleD34
$ws.Range("E34").Value = '  -4.67%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$styleD35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.84'
$ws.Range("D35").Style = $styleD35
$ws.Range("E35").Value = '  +3.04%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$styleD36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0842'
$ws.Range("D36").Style = $styleD36
$ws.Range("E36").Value = '  +8.37%  '

$ws.Range("E37").Value = '  +8.31%  '

$styleD38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '151.98'
$ws.Range("D38").Style = $styleD38
$ws.Range("E38").Value = '  +2.20%  '

$styleD39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.123'
$ws.Range("D39").Style = $styleD39
$ws.Range("E39").Value = '  +7.92%  '

$ws.Range("E40").Value = '  +5.54%  '

$styleD41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.23'
$ws.Range("D41").Style = $styleD41
$ws.Range("E41").Value = '  +39.39%  '

$styleD42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '16.19'
$ws.Range("D42").Style = $styleD42
$ws.Range("E42").Value = '  +5.07%  '

$styleD43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0333'
$ws.Range("D43").Style = $styleD43
$ws.Range("E43").Value = '  +9.76%  '

$styleD44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.63'
$ws.Range("D44").Style = $styleD44
$ws.Range("E44").Value = '  +9.48%  '

$styleD45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.08'
$ws.Range("D45").Style = $styleD45
$ws.Range("E45").Value = '  +3.18%  '

$ws.Range("D46").Value = '2.136.58'
$ws.Range("E46").Value = '  +7.18%  '

$styleD47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.998'
$ws.Range("D47").Style = $styleD47

$styleD48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '93.57'
$ws.Range("D48").Style = $styleD48
$ws.Range("E48").Value = '  +1.01%  '

$styleD49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.57'
$ws.Range("D49").Style = $styleD49
$ws.Range("E49").Value = '  +11.12%  '

$styleD50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.79'
$ws.Range("D50").Style = $styleD50
$ws.Range("E50").Value = '  -1.90%  '

$styleD51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.86'
$ws.Range("D51").Style = $styleD51
$ws.Range("E51").Value = '  +7.51%  '
